$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17,8).Value = 3249
$ws.Cells.Item(17,9).Value = 0
$ws.Cells.Item(17,10).Value = 3249
$ws.Cells.Item(17,11).Value = 0
$ws.Cells.Item(17,12).Value = 9747
$ws.Cells.Item(17,14).Value = -10083
$ws.Cells.Item(64,8).Value = 32865.816
$ws.Cells.Item(64,9).Value = 57859.39
$ws.Cells.Item(64,10).Value = 2873.5334
$ws.Cells.Item(64,11).Value = 57859.39
$ws.Cells.Item(64,12).Value = 2873.5334
$ws.Cells.Item(64,13).Value = -57611.39
$ws.Cells.Item(64,14).Value = -3369.5334
$ws.Cells.Item(67,8).Value = 32865.816
$ws.Cells.Item(67,9).Value = 57859.39
$ws.Cells.Item(67,10).Value = 2873.5334
$ws.Cells.Item(67,11).Value = 57859.39
$ws.Cells.Item(67,12).Value = 2873.5334
$ws.Cells.Item(67,13).Value = -57001.39
$ws.Cells.Item(67,14).Value = -4589.5334
$ws.Cells.Item(70,8).Value = 1500
$ws.Cells.Item(70,9).Value = 0
$ws.Cells.Item(70,10).Value = 1500
$ws.Cells.Item(70,11).Value = 0
$ws.Cells.Item(70,12).Value = 4500
$ws.Cells.Item(70,13).Value = ""
$ws.Cells.Item(70,14).Value = -5040
$ws.Cells.Item(73,8).Value = 1500
$ws.Cells.Item(73,9).Value = 0
$ws.Cells.Item(73,10).Value = 1500
$ws.Cells.Item(73,11).Value = 0
$ws.Cells.Item(73,12).Value = 4500
$ws.Cells.Item(73,13).Value = ""
$ws.Cells.Item(73,14).Value = -6372
$ws.Cells.Item(103,8).Value = 1071.1428
$ws.Cells.Item(103,9).Value = 1354.5
$ws.Cells.Item(103,10).Value = 693.3333
$ws.Cells.Item(103,11).Value = 4063.5
$ws.Cells.Item(103,12).Value = 2079.9999
$ws.Cells.Item(103,13).Value = -3477.5
$ws.Cells.Item(103,14).Value = -3251.9999
$ws.Cells.Item(112,8).Value = 1311.2
$ws.Cells.Item(112,9).Value = 850
$ws.Cells.Item(112,10).Value = 1344.1428
$ws.Cells.Item(112,11).Value = 2550
$ws.Cells.Item(112,12).Value = 4032.4284
$ws.Cells.Item(112,13).Value = -1442
$ws.Cells.Item(112,14).Value = -6248.428400000001
$ws.Cells.Item(132,8).Value = 15920.954
$ws.Cells.Item(132,9).Value = 2336.6606
$ws.Cells.Item(132,10).Value = 100445.445
$ws.Cells.Item(132,11).Value = 7009.9818
$ws.Cells.Item(132,12).Value = 301336.335
$ws.Cells.Item(132,13).Value = -4479.9818
$ws.Cells.Item(132,14).Value = -306396.335
$ws.Cells.Item(138,8).Value = 1791.3738
$ws.Cells.Item(138,9).Value = 1478.5834
$ws.Cells.Item(138,10).Value = 1970.1111
$ws.Cells.Item(138,11).Value = 4435.7502
$ws.Cells.Item(138,12).Value = 5910.3333
$ws.Cells.Item(138,13).Value = 704.2497999999996
$ws.Cells.Item(138,14).Value = -16190.3333
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 31488.285
$ws.Cells.Item(32,9).Value = 33627.414
$ws.Cells.Item(32,10).Value = 21648.3
$ws.Cells.Item(32,11).Value = 33627.414
$ws.Cells.Item(32,12).Value = 21648.3
$ws.Cells.Item(32,13).Value = -33340.414
$ws.Cells.Item(32,14).Value = -22222.3
$ws.Cells.Item(61,8).Value = 3459.2
$ws.Cells.Item(61,9).Value = 2275.5386
$ws.Cells.Item(61,10).Value = 4364.353
$ws.Cells.Item(61,11).Value = 2275.5386
$ws.Cells.Item(61,12).Value = 4364.353
$ws.Cells.Item(61,13).Value = -2063.5386
$ws.Cells.Item(61,14).Value = -4788.353
$ws.Cells.Item(102,8).Value = 13939
$ws.Cells.Item(102,9).Value = 1831.5385
$ws.Cells.Item(102,10).Value = 53288.25
$ws.Cells.Item(102,11).Value = 1831.5385
$ws.Cells.Item(102,12).Value = 53288.25
$ws.Cells.Item(102,13).Value = -209.5385000000001
$ws.Cells.Item(102,14).Value = -56532.25
$ws.Cells.Item(121,8).Value = 29656.625
$ws.Cells.Item(121,9).Value = 0
$ws.Cells.Item(121,10).Value = 29656.625
$ws.Cells.Item(121,11).Value = 0
$ws.Cells.Item(121,12).Value = 29656.625
$ws.Cells.Item(121,14).Value = -33150.625
$ws.Cells.Item(122,8).Value = 1787.7826
$ws.Cells.Item(122,9).Value = 1766.8
$ws.Cells.Item(122,10).Value = 1854.5454
$ws.Cells.Item(122,11).Value = 5300.4
$ws.Cells.Item(122,12).Value = 5563.6362
$ws.Cells.Item(122,13).Value = -2850.4
$ws.Cells.Item(122,14).Value = -10463.6362
$ws.Cells.Item(136,8).Value = 3459.2
$ws.Cells.Item(136,9).Value = 2275.5386
$ws.Cells.Item(136,10).Value = 4364.353
$ws.Cells.Item(136,11).Value = 6826.6158
$ws.Cells.Item(136,12).Value = 13093.059
$ws.Cells.Item(136,13).Value = -4276.6158
$ws.Cells.Item(136,14).Value = -18193.059
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107,8).Value = 2065.3845
$ws.Cells.Item(107,9).Value = 1968.8462
$ws.Cells.Item(107,10).Value = 2161.923
$ws.Cells.Item(107,11).Value = 1968.8462
$ws.Cells.Item(107,12).Value = 2161.923
$ws.Cells.Item(107,13).Value = -48.84619999999995
$ws.Cells.Item(107,14).Value = -6001.923
$ws.Cells.Item(112,8).Value = 46361.25
$ws.Cells.Item(112,9).Value = 0
$ws.Cells.Item(112,10).Value = 46361.25
$ws.Cells.Item(112,11).Value = 0
$ws.Cells.Item(112,12).Value = 46361.25
$ws.Cells.Item(112,14).Value = -49315.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 159816.42
$ws.Cells.Item(31,9).Value = 1855.6522
$ws.Cells.Item(31,10).Value = 214041.77
$ws.Cells.Item(31,11).Value = 1855.6522
$ws.Cells.Item(31,12).Value = 214041.77
$ws.Cells.Item(31,13).Value = -1560.6522
$ws.Cells.Item(31,14).Value = -214631.77
$ws.Cells.Item(34,8).Value = 159816.42
$ws.Cells.Item(34,9).Value = 1855.6522
$ws.Cells.Item(34,10).Value = 214041.77
$ws.Cells.Item(34,11).Value = 1855.6522
$ws.Cells.Item(34,12).Value = 214041.77
$ws.Cells.Item(34,13).Value = -1653.6522
$ws.Cells.Item(34,14).Value = -214445.77
$ws.Cells.Item(99,8).Value = 1722.4706
$ws.Cells.Item(99,9).Value = 1410.25
$ws.Cells.Item(99,10).Value = 2000
$ws.Cells.Item(99,11).Value = 1410.25
$ws.Cells.Item(99,12).Value = 2000
$ws.Cells.Item(99,13).Value = 87.75
$ws.Cells.Item(99,14).Value = -4996
$ws.Cells.Item(126,8).Value = 1722.4706
$ws.Cells.Item(126,9).Value = 1410.25
$ws.Cells.Item(126,10).Value = 2000
$ws.Cells.Item(126,11).Value = 4230.75
$ws.Cells.Item(126,12).Value = 6000
$ws.Cells.Item(126,13).Value = -1760.75
$ws.Cells.Item(126,14).Value = -10940
$ws.Cells.Item(132,8).Value = 54775.223
$ws.Cells.Item(132,9).Value = 1790.4
$ws.Cells.Item(132,10).Value = 206160.42
$ws.Cells.Item(132,11).Value = 5371.200000000001
$ws.Cells.Item(132,12).Value = 618481.26
$ws.Cells.Item(132,13).Value = -2841.200000000001
$ws.Cells.Item(132,14).Value = -623541.26
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102,8).Value = 2748.1428
$ws.Cells.Item(102,9).Value = 2444.8
$ws.Cells.Item(102,10).Value = 3506.5
$ws.Cells.Item(102,11).Value = 2444.8
$ws.Cells.Item(102,12).Value = 3506.5
$ws.Cells.Item(102,13).Value = -822.8000000000002
$ws.Cells.Item(102,14).Value = -6750.5
$ws.Cells.Item(104,8).Value = 33393.8
$ws.Cells.Item(104,9).Value = 0
$ws.Cells.Item(104,10).Value = 33393.8
$ws.Cells.Item(104,11).Value = 0
$ws.Cells.Item(104,12).Value = 33393.8
$ws.Cells.Item(104,14).Value = -40381.8
$ws.Cells.Item(110,8).Value = 31075.2
$ws.Cells.Item(110,9).Value = 0
$ws.Cells.Item(110,10).Value = 31075.2
$ws.Cells.Item(110,11).Value = 0
$ws.Cells.Item(110,12).Value = 31075.2
$ws.Cells.Item(110,14).Value = -39255.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82,8).Value = 6945500
$ws.Cells.Item(82,9).Value = 1196.1111
$ws.Cells.Item(82,10).Value = 27778412
$ws.Cells.Item(82,11).Value = 1196.1111
$ws.Cells.Item(82,12).Value = 27778412
$ws.Cells.Item(82,13).Value = -835.1111000000001
$ws.Cells.Item(82,14).Value = -27779134
$ws.Cells.Item(85,8).Value = 6945500
$ws.Cells.Item(85,9).Value = 1196.1111
$ws.Cells.Item(85,10).Value = 27778412
$ws.Cells.Item(85,11).Value = 1196.1111
$ws.Cells.Item(85,12).Value = 27778412
$ws.Cells.Item(85,13).Value = 51.88889999999992
$ws.Cells.Item(85,14).Value = -27780908
$ws.Cells.Item(106,8).Value = 30643.334
$ws.Cells.Item(106,9).Value = 0
$ws.Cells.Item(106,10).Value = 30643.334
$ws.Cells.Item(106,11).Value = 0
$ws.Cells.Item(106,12).Value = 30643.334
$ws.Cells.Item(106,14).Value = -33167.334
$ws.Cells.Item(121,8).Value = 18891.666
$ws.Cells.Item(121,9).Value = 0
$ws.Cells.Item(121,10).Value = 18891.666
$ws.Cells.Item(121,11).Value = 0
$ws.Cells.Item(121,12).Value = 18891.666
$ws.Cells.Item(121,14).Value = -22385.666
$ws.Cells.Item(132,8).Value = 4161.2383
$ws.Cells.Item(132,9).Value = 2488.5557
$ws.Cells.Item(132,10).Value = 5415.75
$ws.Cells.Item(132,11).Value = 7465.6671
$ws.Cells.Item(132,12).Value = 16247.25
$ws.Cells.Item(132,13).Value = -4935.6671
$ws.Cells.Item(132,14).Value = -21307.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16,8).Value = 46804
$ws.Cells.Item(16,9).Value = 0
$ws.Cells.Item(16,10).Value = 46804
$ws.Cells.Item(16,11).Value = 0
$ws.Cells.Item(16,12).Value = 46804
$ws.Cells.Item(16,14).Value = -47388
$ws.Cells.Item(96,8).Value = 1842.125
$ws.Cells.Item(96,9).Value = 990
$ws.Cells.Item(96,10).Value = 2694.25
$ws.Cells.Item(96,11).Value = 990
$ws.Cells.Item(96,12).Value = 2694.25
$ws.Cells.Item(96,13).Value = 383
$ws.Cells.Item(96,14).Value = -5440.25
$ws.Cells.Item(121,8).Value = 44412
$ws.Cells.Item(121,9).Value = 0
$ws.Cells.Item(121,10).Value = 44412
$ws.Cells.Item(121,11).Value = 0
$ws.Cells.Item(121,12).Value = 44412
$ws.Cells.Item(121,14).Value = -47906
$ws.Cells.Item(122,8).Value = 818.9231
$ws.Cells.Item(122,9).Value = 764.7
$ws.Cells.Item(122,10).Value = 999.6667
$ws.Cells.Item(122,11).Value = 2294.1
$ws.Cells.Item(122,12).Value = 2999.0001
$ws.Cells.Item(122,13).Value = 155.8999999999996
$ws.Cells.Item(122,14).Value = -7899.0001
$ws.Cells.Item(126,8).Value = 2534.818
$ws.Cells.Item(126,9).Value = 2355.842
$ws.Cells.Item(126,10).Value = 3668.3333
$ws.Cells.Item(126,11).Value = 7067.526
$ws.Cells.Item(126,12).Value = 11004.9999
$ws.Cells.Item(126,13).Value = -4597.526
$ws.Cells.Item(126,14).Value = -15944.9999
$ws.Cells.Item(131,8).Value = 54334
$ws.Cells.Item(131,9).Value = 0
$ws.Cells.Item(131,10).Value = 54334
$ws.Cells.Item(131,11).Value = 0
$ws.Cells.Item(131,12).Value = 54334
$ws.Cells.Item(131,14).Value = -64414
